$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - column headers for the job posting sheet
$ws.Range("A1").Value = "Job_Id"
$ws.Range("B1").Value = "Job_Title"
$ws.Range("C1").Value = "Job_Description"
$ws.Range("D1").Value = "Total_Years_Min_Exp"
$ws.Range("E1").Value = "Total_Years_Max_Exp"
$ws.Range("F1").Value = "LinkedIn_Poster"
$ws.Range("G1").Value = "LinkedIn_Posted"
$ws.Range("H1").Value = "Resume_received"
$ws.Range("I1").Value = "Resume_downloaded"

# Data row (row 2) - the new job posting JD_001
$ws.Range("A2").Value = "JD_001"
$ws.Range("B2").Value = "Senior Dotnet Engineer"
$ws.Range("C2").Value = "A senior .NET developer is responsible for designing, developing, and maintaining high-quality, scalable applications on the .NET framework, from conception to deployment. Key duties include collaborating with cross-functional teams, writing efficient code, mentoring junior developers, and ensuring high performance and security. This role involves a mix of technical leadership and hands-on development, requiring proficiency in core .NET technologies, architecture, databases, and agile methodologies. "
$ws.Range("D2").Value = 2
$ws.Range("E2").Value = 4
